$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "SELECT`n    COUNT(DISTINCT p.program_acronym) AS `"Programs`",`n    COUNT(DISTINCT st.clinical_study_designation) AS `"Studies`",`n    COUNT(DISTINCT c.case_record_id) AS `"Cases`",`n    COUNT(DISTINCT smp.sample_id) AS `"Samples`",`n    COUNT(DISTINCT cf.file_name) AS `"Case Files`",         `n    COUNT(DISTINCT sf.file_name) AS `"Study Files`" `nFROM `n    df_program p`nJOIN `n    df_study st ON p.program_acronym = st.`"program.program_acronym`"`nJOIN `n    df_case c ON st.clinical_study_designation = c.`"study.clinical_study_designation`"`nJOIN `n    df_demographic dmg ON dmg.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_diagnosis diag ON diag.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_enrollment enr ON enr.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_sample smp ON smp.`"case.case_record_id`" = c.case_record_id`nLEFT JOIN `n    df_case_file cf ON cf.`"sample.sample_id`" = smp.sample_id`nLEFT JOIN `n    df_study_file sf ON sf.`"study.clinical_study_designation`" = st.clinical_study_designation`nWHERE `n    st.clinical_study_designation = 'MGT01' AND cf.file_type = 'Index File';"
$ws.Range("B2").Value = "SELECT`n    DISTINCT c.case_record_id AS `"Case ID`",`n    st.clinical_study_designation AS `"Study Code`",`n    st.clinical_study_type AS `"Study Type`",`n    dmg.breed AS `"Breed`",`n    diag.disease_term AS `"Diagnosis`",`n    diag.stage_of_disease AS `"Stage Of Disease`",`n    CASE `n    WHEN dmg.patient_age_at_enrollment = CAST(dmg.patient_age_at_enrollment AS INT) `n    THEN CAST(CAST(dmg.patient_age_at_enrollment AS INT) AS TEXT)`n    ELSE CAST(dmg.patient_age_at_enrollment AS TEXT)`nEND AS `"Age`", `n    COALESCE(TRIM(dmg.sex), '') AS `"Sex`",   `n    COALESCE(dmg.neutered_indicator, '') AS `"Neutered Status`",  `n    COALESCE(`n        CASE `n            WHEN dmg.weight = CAST(dmg.weight AS INT) THEN CAST(dmg.weight AS INT)`n            ELSE dmg.weight `n        END, `n    '') AS `"Weight (kg)`",        `n    COALESCE(diag.best_response, '') AS `"Response to Treatment`",  `n    COALESCE(coh.cohort_description, '') AS `"Cohort`"  `nFROM `n    df_case c`nJOIN `n    df_study st ON c.`"study.clinical_study_designation`" = st.clinical_study_designation`nJOIN `n    df_demographic dmg ON dmg.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_diagnosis diag ON diag.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_enrollment enr ON enr.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_program p ON st.`"program.program_acronym`" = p.program_acronym`nJOIN `n    df_sample smp ON smp.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_publication pub ON pub.`"study.clinical_study_designation`" = st.clinical_study_designation`nLEFT JOIN `n    df_case_file cf ON cf.`"sample.sample_id`" = smp.sample_id`nLEFT JOIN `n    df_study_file sf ON sf.`"study.clinical_study_designation`" = st.clinical_study_designation`nLEFT JOIN`n    df_cohort coh ON coh.`"study.clinical_study_designation`" = st.clinical_study_designation`nWHERE`n    st.clinical_study_designation = 'MGT01' AND cf.file_type = 'Index File'`nORDER BY `n    c.case_record_id ASC`nLIMIT 100;"
$ws.Range("B3").Value = "SELECT DISTINCT`n    smp.sample_id AS `"Sample ID`",`n    c.case_record_id AS `"Case ID`",`n    COALESCE(dmg.breed, '') AS `"Breed`",`n    COALESCE(diag.disease_term, '') AS `"Diagnosis`",`n    COALESCE(smp.sample_site, '') AS `"Sample Site`",`n    COALESCE(smp.summarized_sample_type, '') AS `"Sample Type`",`n    COALESCE(smp.specific_sample_pathology, '') AS `"Pathology/Morphology`",`n    COALESCE(smp.tumor_grade, '') AS `"Tumor Grade`",`n    COALESCE(smp.sample_chronology, '') AS `"Sample Chronology`",`n    COALESCE(smp.percentage_tumor, '') AS `"Percentage Tumor`",`n    COALESCE(smp.necropsy_sample, '') AS `"Necropsy Sample`",`n    COALESCE(smp.sample_preservation, '') AS `"Sample Preservation`"`nFROM `n    df_sample smp`nJOIN `n    df_case c ON smp.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_publication pub ON pub.`"study.clinical_study_designation`" = st.clinical_study_designation`nJOIN `n    df_demographic dmg ON dmg.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_diagnosis diag ON diag.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_enrollment enr ON enr.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_program p ON st.`"program.program_acronym`" = p.program_acronym`nJOIN `n    df_study st ON c.`"study.clinical_study_designation`" = st.clinical_study_designation`nLEFT JOIN `n    df_case_file cf ON cf.`"sample.sample_id`" = smp.sample_id`nLEFT JOIN `n    df_study_file sf ON sf.`"study.clinical_study_designation`" = st.clinical_study_designation`nWHERE `n   st.clinical_study_designation = 'MGT01' AND cf.file_type = 'Index File'`nORDER BY `n    smp.sample_id ASC`nLIMIT 100;"
$ws.Range("B4").Value = "SELECT `n    DISTINCT cf.file_name AS `"File Name`",`n    CASE`n        WHEN cf.file_name LIKE '%.bai' THEN 'bai'`n        WHEN cf.file_name LIKE '%.bam' THEN 'bam'`n        WHEN cf.file_name LIKE '%.csv' THEN 'csv'`n        WHEN cf.file_name LIKE '%.doc' THEN 'doc'`n        WHEN cf.file_name LIKE '%.docx' THEN 'docx'`n        WHEN cf.file_name LIKE '%.gz' THEN 'gz'`n        WHEN cf.file_name LIKE '%.pdf' THEN 'pdf'`n        WHEN cf.file_name LIKE '%.rtf' THEN 'rtf'`n        WHEN cf.file_name LIKE '%.tbi' THEN 'tbi'`n        WHEN cf.file_name LIKE '%.tif' THEN 'tif'`n        WHEN cf.file_name LIKE '%.xls' THEN 'xls'`n        WHEN cf.file_name LIKE '%.xlsx' THEN 'xlsx'`n        ELSE 'Unknown'`n    END AS `"Format`",`n    cf.file_type AS `"File Type`",`n    CASE     `n    WHEN cf.file_size >= 1024 * 1024 * 1024 THEN `n        CASE `n            WHEN ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) = CAST(ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) `n            THEN CAST(CAST(ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' GB'`n            ELSE ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB'`n        END`n    WHEN cf.file_size >= 1024 * 1024 THEN `n        CASE `n            WHEN ROUND(cf.file_size / (1024.0 * 1024.0), 2) = CAST(ROUND(cf.file_size / (1024.0 * 1024.0), 0) AS INT) `n            THEN CAST(CAST(ROUND(cf.file_size / (1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' MB'`n            ELSE ROUND(cf.file_size / (1024.0 * 1024.0), 2) || ' MB'`n        END`n    WHEN cf.file_size >= 1024 THEN `n        CASE `n            WHEN ROUND(cf.file_size / 1024.0, 2) = CAST(ROUND(cf.file_size / 1024.0, 0) AS INT) `n            THEN CAST(CAST(ROUND(cf.file_size / 1024.0, 0) AS INT) AS TEXT) || ' KB'`n            ELSE ROUND(cf.file_size / 1024.0, 2) || ' KB'`n        END`n    ELSE `n        CASE `n            WHEN ROUND(cf.file_size, 2) = CAST(ROUND(cf.file_size, 0) AS INT) `n            THEN CAST(CAST(ROUND(cf.file_size, 0) AS INT) AS TEXT) || ' Bytes'`n            ELSE ROUND(cf.file_size, 2) || ' Bytes'`n        END`nEND AS `"Size`",`n    'sample' AS `"Association`",`n    cf.file_description AS `"Description`",`n    smp.sample_id AS `"Sample ID`",`n    c.case_record_id AS `"Case ID`",`n    dmg.breed AS `"Breed`",`n    diag.disease_term AS `"Diagnosis`"`nFROM `n    df_case_file cf`nJOIN `n    df_sample smp ON cf.`"sample.sample_id`" = smp.sample_id`nJOIN `n    df_case c ON smp.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_study st ON c.`"study.clinical_study_designation`" = st.clinical_study_designation`nJOIN `n    df_program p ON st.`"program.program_acronym`" = p.program_acronym`nJOIN `n    df_demographic dmg ON dmg.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_diagnosis diag ON diag.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_enrollment enr ON enr.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_publication pub ON pub.`"study.clinical_study_designation`" = st.clinical_study_designation`nLEFT JOIN `n    df_study_file sf ON sf.`"study.clinical_study_designation`" = st.clinical_study_designation`nWHERE`n    st.clinical_study_designation = 'MGT01' AND cf.file_type = 'Index File'`nORDER BY `n    cf.file_name ASC`nLIMIT 100;"
$ws.Range("B5").Value = "SELECT DISTINCT`n    sf.file_name AS `"File Name`",`n    sf.file_type AS `"File Type`",`n    'study' AS `"Association`",`n    sf.file_description AS `"Description`",`n    CASE`n        WHEN sf.file_name LIKE '%.bai' THEN 'bai'`n        WHEN sf.file_name LIKE '%.bam' THEN 'bam'`n        WHEN sf.file_name LIKE '%.csv' THEN 'csv'`n        WHEN sf.file_name LIKE '%.doc' THEN 'doc'`n        WHEN sf.file_name LIKE '%.docx' THEN 'docx'`n        WHEN sf.file_name LIKE '%.gz' THEN 'gz'`n        WHEN sf.file_name LIKE '%.pdf' THEN 'pdf'`n        WHEN sf.file_name LIKE '%.rtf' THEN 'rtf'`n        WHEN sf.file_name LIKE '%.tbi' THEN 'tbi'`n        WHEN sf.file_name LIKE '%.tif' THEN 'tif'`n        WHEN sf.file_name LIKE '%.xls' THEN 'xls'`n        WHEN sf.file_name LIKE '%.xlsx' THEN 'xlsx'`n        ELSE 'Unknown'`n    END AS `"Format`",`n    CASE     `n        WHEN sf.file_size >= 1024 * 1024 * 1024 THEN `n            ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB' `n        WHEN sf.file_size >= 1024 * 1024 THEN `n            ROUND(sf.file_size / (1024.0 * 1024.0), 2) || ' MB' `n        WHEN sf.file_size >= 1024 THEN `n            ROUND(sf.file_size / 1024.0, 2) || ' KB' `n        ELSE `n            ROUND(sf.file_size, 2) || ' Bytes' `n    END AS `"Size`",`n    st.clinical_study_designation AS `"Study Code`"`nFROM `n    df_case_file cf`nJOIN `n    df_sample smp ON cf.`"sample.sample_id`" = smp.sample_id`nJOIN `n    df_case c ON smp.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_study st ON c.`"study.clinical_study_designation`" = st.clinical_study_designation`nJOIN `n    df_program p ON st.`"program.program_acronym`" = p.program_acronym`nJOIN `n    df_demographic dmg ON dmg.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_diagnosis diag ON diag.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_enrollment enr ON enr.`"case.case_record_id`" = c.case_record_id`nJOIN `n    df_publication pub ON pub.`"study.clinical_study_designation`" = st.clinical_study_designation`nLEFT JOIN `n    df_study_file sf ON sf.`"study.clinical_study_designation`" = st.clinical_study_designation`nWHERE`n    st.clinical_study_designation = 'MGT01' AND cf.file_type = 'Index File'`nORDER BY `n    sf.file_name ASC`nLIMIT 100;"
